$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GVA")

$ws.Range("B2").Value = 0.02927166946637607
$ws.Range("C2").Value = 1.520652503372878
$ws.Range("D2").Value = 6.746610706108599
$ws.Range("E2").Value = 2.597423859540179
$ws.Range("F2").Value = 2.658379252570611
$ws.Range("G2").Value = 22
